$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.671.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "'1.946.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'247.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4806"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").Value = "'0.2921"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "'0.06796"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'111.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "'19.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "'1.927.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "'0.07672"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "'5.486"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.87%  "
$ws.Range("D15").Value = "'0.6849"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "'292.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").Value = "'30.664.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "'13.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").Value = "'5.641"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").Value = "'0.000007665"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "'2.203.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'6.571"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").Value = "'9.754"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").Value = "'168.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").Value = "'20.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "'2.177"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").Value = "'1.433"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("D31").Value = "'4.689"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.81%  "
$ws.Range("D32").Value = "'4.491"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.06%  "
$ws.Range("D33").Value = "'0.05047"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'0.7687"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("D35").Value = "'1.155"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("D36").Value = "'0.02066"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").Value = "'2.735"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'2.697"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").Value = "'2.044"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Value = "'110.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "'0.4450"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'0.8690"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "'5.938"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").Value = "'69.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").Value = "'7.340"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'9.383"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "'48.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'0.1249"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "'35.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").Value = "'0.2508"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.69%  "
